# The author removed the post "「私できるよ」" which occupied row 784,
# causing every subsequent row (785-798) to shift up by one and the
# sheet's used range to shrink from A1:C798 to A1:C797.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("784").Delete()
